$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 0.1.1 -> 0.1.2
$ws.Range("B3").Value = "0.1.2"

# Title: "GVHD ICD-10 Codes" -> "ICD-10 (GVHD) Codes"
$ws.Range("B5").Value = "ICD-10 (GVHD) Codes"

# Date: 2024-12-02T18:31:42-06:00 -> 2025-04-16T10:37:17-05:00
$ws.Range("B8").Value = "2025-04-16T10:37:17-05:00"
